$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value that was updated from
# 45192 (2023-09-23) to 45202 (2023-10-03) for every data row (rows 2-79).
$ws.Range("C2:C79").Value = 45202
